# Update the LinkedIn carousel draft: replace the Juniper Green Energy
# article text with the new NTPC Green Energy / Bhadla Solar Project copy
# across all six slides, preserving existing run/paragraph formatting by
# only rewriting the substring of each shape's TextRange that changed.

function Replace-ShapeText {
    param($Shape, $OldText, $NewText)

    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Text not found: $OldText"
    }
    $start = $idx + 1
    $len = $OldText.Length
    $tr.Characters($start, $len).Text = $NewText
}

$p = $ppt.ActivePresentation

$oldTitle = "Juniper Green Energy Commissions Additional 72 MWp Solar Component of Hybrid Project in Maharashtra - Energetica India Magazine"
$newTitle = "NTPC Green Energy Commissions 300 MW Phase Of 500 MW Bhadla Solar Project In Rajasthan - SolarQuarter"

$bodyReplacements = @{
    1 = @(
        @("Juniper Green Energy has commissioned an additional 72 MWp solar component.", "NTPC Green Energy has commissioned the 300 MW phase of the Bhadla Solar Project."),
        @("The project is part of a hybrid initiative located in Maharashtra.", "The total capacity of the Bhadla Solar Project is 500 MW.")
    )
    2 = @(
        @("The total capacity of the solar component is 72 MWp.", "The Bhadla Solar Project is located in Rajasthan, India."),
        @("This addition enhances the renewable energy output in Maharashtra.", "This commissioning is part of NTPC's efforts to expand its renewable energy portfolio.")
    )
    3 = @(
        @("The hybrid project aims to integrate solar energy with other renewable sources.", "The 300 MW phase is a significant step towards achieving NTPC's renewable energy targets."),
        @("Maharashtra is a key region for renewable energy development in India.", "NTPC aims to increase its renewable energy capacity in the coming years.")
    )
    4 = @(
        @("The commissioning of the solar component contributes to India's renewable energy targets.", "The commissioning of this phase contributes to India's solar energy generation capacity."),
        @("Juniper Green Energy is focused on expanding its renewable energy portfolio.", "Bhadla Solar Park is one of the largest solar parks in the world.")
    )
    5 = @(
        @("The project reflects ongoing investments in sustainable energy solutions.", "This project aligns with India's commitment to increase renewable energy usage."),
        @("This initiative supports local energy needs and reduces carbon footprint.", "NTPC is a key player in India's transition to sustainable energy sources.")
    )
    6 = @(
        @("The hybrid project is part of a broader strategy to enhance energy security.", "The project is part of NTPC's broader strategy to enhance energy security."),
        @("Juniper Green Energy's efforts align with national policies on renewable energy.", "The completion of this phase marks a milestone in NTPC's solar initiatives.")
    )
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = $slide.Shapes.Item(1)
    $bodyShape = $slide.Shapes.Item(2)

    Replace-ShapeText $titleShape $oldTitle $newTitle

    foreach ($pair in $bodyReplacements[$i]) {
        Replace-ShapeText $bodyShape $pair[0] $pair[1]
    }
}

Write-Host "Done updating $($p.Slides.Count) slides"
